$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new exam papers (9618_s24_42 and 9618_s24_43) need to be tracked at
# the top of the list. Make room for them by nudging the existing
# "Question Paper" / "Q1" / "Q2" / "Q3" entries (columns B:E) down two
# rows (current rows 6-11 -> rows 8-13) - without touching the S.No.
# column (A) or any row below the shifted block.

# Snapshot the block that is about to move, bottom-up isn't needed since
# we read everything into memory before writing anything back.
$snapshot = @()
for ($row = 6; $row -le 11; $row++) {
    $snapshot += , @(
        $ws.Cells.Item($row, 2).Value(),
        $ws.Cells.Item($row, 3).Value(),
        $ws.Cells.Item($row, 4).Value(),
        $ws.Cells.Item($row, 5).Value()
    )
}

# Wipe the block that is about to be rewritten so no stray values are
# left behind once the data has shifted down.
$ws.Range("B6:E13").ClearContents()

for ($i = 0; $i -lt $snapshot.Count; $i++) {
    $row = $i + 8
    $data = $snapshot[$i]
    if ($data[0]) { $ws.Cells.Item($row, 2).Value = $data[0] }
    if ($data[1]) { $ws.Cells.Item($row, 3).Value = $data[1] }
    if ($data[2]) { $ws.Cells.Item($row, 4).Value = $data[2] }
    if ($data[3]) { $ws.Cells.Item($row, 5).Value = $data[3] }
}

# Drop the two new question papers into the freed-up rows.
$ws.Range("B7").Value = "9618_s24_43"
$ws.Range("B6").Value = "9618_s24_42"

# Selection ends up on A16 after the edit.
$ws.Range("A16").Select()
